$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (ion_id / vector-string columns): A, C, I, K ---
# Force Text number format on the touched cells first so values round-trip as
# strings rather than being auto-coerced to numbers, then clear the format again
# so the cell style matches the original (unstyled) cells.
$textCells = @(
    "A2",
    "C2",
    "I2",
    "K2",
    "A3",
    "C3",
    "I3",
    "K3",
    "A4",
    "C4",
    "I4",
    "K4",
    "A5",
    "C5",
    "I5",
    "K5",
    "K6",
    "A7",
    "C7",
    "I7",
    "K7",
    "A8",
    "C8",
    "I8",
    "K8",
    "A9",
    "C9",
    "I9",
    "K9",
    "A10",
    "C10",
    "I10",
    "K10"
)
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("A2").Value = "1307"
$ws.Range("C2").Value = "[-0.5532672293484211, 0.28407424688339233, -8.042064785957336]"
$ws.Range("I2").Value = "[0.6336288452148438, -1.4450759887695312, -3.3725357055664062]"
$ws.Range("K2").Value = "[-0.008810419589281082, -2.18619704246521, -2.864572048187256]"
$ws.Range("A3").Value = "1307"
$ws.Range("C3").Value = "[-0.5532672293484211, 0.28407424688339233, -8.042064785957336]"
$ws.Range("I3").Value = "[0.6336288452148438, -1.4450759887695312, -3.3725357055664062]"
$ws.Range("K3").Value = "[-0.2788090407848358, 0.07740062475204468, -1.6099847555160522]"
$ws.Range("A4").Value = "1307"
$ws.Range("C4").Value = "[-0.5532672293484211, 0.28407424688339233, -8.042064785957336]"
$ws.Range("I4").Value = "[0.6336288452148438, -1.4450759887695312, -3.3725357055664062]"
$ws.Range("K4").Value = "[0.6154102087020874, 1.285620093345642, -2.1181087493896484]"
$ws.Range("A5").Value = "1307"
$ws.Range("C5").Value = "[-0.5532672293484211, 0.28407424688339233, -8.042064785957336]"
$ws.Range("I5").Value = "[0.6336288452148438, -1.4450759887695312, -3.3725357055664062]"
$ws.Range("K5").Value = "[-0.8810579776763916, 1.1072505712509155, -1.4493992328643799]"
$ws.Range("K6").Value = "[-0.05633280798792839, 0.0005689358222298324, -1.5595158338546753]"
$ws.Range("A7").Value = "1306"
$ws.Range("C7").Value = "[2.1177535615861416, 1.4151039516436867, -12.508464455604553]"
$ws.Range("I7").Value = "[-0.22468185424804688, 0.5847549438476562, -10.533576965332031]"
$ws.Range("K7").Value = "[0.07805076986551285, -0.4003472626209259, -3.418886423110962]"
$ws.Range("A8").Value = "1306"
$ws.Range("C8").Value = "[2.1177535615861416, 1.4151039516436867, -12.508464455604553]"
$ws.Range("I8").Value = "[-0.22468185424804688, 0.5847549438476562, -10.533576965332031]"
$ws.Range("K8").Value = "[-1.893072247505188, -3.221282482147217, -2.009432315826416]"
$ws.Range("A9").Value = "1306"
$ws.Range("C9").Value = "[2.1177535615861416, 1.4151039516436867, -12.508464455604553]"
$ws.Range("I9").Value = "[-0.22468185424804688, 0.5847549438476562, -10.533576965332031]"
$ws.Range("K9").Value = "[1.36846923828125, 8.766291618347168, -3.4916884899139404]"
$ws.Range("A10").Value = "1306"
$ws.Range("C10").Value = "[2.1177535615861416, 1.4151039516436867, -12.508464455604553]"
$ws.Range("I10").Value = "[-0.22468185424804688, 0.5847549438476562, -10.533576965332031]"
$ws.Range("K10").Value = "[2.620638608932495, -3.7301268577575684, -2.0289413928985596]"

foreach ($addr in $textCells) { $ws.Range($addr).ClearFormats() }

# --- Numeric-valued cells ---
$ws.Range("B2").Value = 5178
$ws.Range("D2").Value = 8.066077660570116
$ws.Range("E2").Value = 7.079835093630851
$ws.Range("F2").Value = 0.8777295969067627
$ws.Range("G2").Value = 0.6219347271323936
$ws.Range("H2").Value = -8.042064785957336
$ws.Range("J2").Value = 1320
$ws.Range("L2").Value = 3.603513479232788
$ws.Range("M2").Value = 9.598557472229004
$ws.Range("N2").Value = 0.9550714863493811
$ws.Range("O2").Value = 3.441612809865843
$ws.Range("P2").Value = 98
$ws.Range("Q2").Value = 98
$ws.Range("R2").Value = 98
$ws.Range("B3").Value = 5178
$ws.Range("D3").Value = 8.066077660570116
$ws.Range("E3").Value = 7.079835093630851
$ws.Range("F3").Value = 0.8777295969067627
$ws.Range("G3").Value = 0.6219347271323936
$ws.Range("H3").Value = -8.042064785957336
$ws.Range("J3").Value = 1309
$ws.Range("L3").Value = 1.635780096054077
$ws.Range("M3").Value = 14.24644756317139
$ws.Range("N3").Value = 0.8441143479128079
$ws.Range("O3").Value = 1.380785379002085
$ws.Range("B4").Value = 5178
$ws.Range("D4").Value = 8.066077660570116
$ws.Range("E4").Value = 7.079835093630851
$ws.Range("F4").Value = 0.8777295969067627
$ws.Range("G4").Value = 0.6219347271323936
$ws.Range("H4").Value = -8.042064785957336
$ws.Range("J4").Value = 1308
$ws.Range("L4").Value = 2.553024530410767
$ws.Range("M4").Value = 11.40358924865723
$ws.Range("N4").Value = 0.597050028512469
$ws.Range("O4").Value = 1.524283271710701
$ws.Range("P4").Value = 1073
$ws.Range("Q4").Value = 1073
$ws.Range("R4").Value = 1073
$ws.Range("B5").Value = 5178
$ws.Range("D5").Value = 8.066077660570116
$ws.Range("E5").Value = 7.079835093630851
$ws.Range("F5").Value = 0.8777295969067627
$ws.Range("G5").Value = 0.6219347271323936
$ws.Range("H5").Value = -8.042064785957336
$ws.Range("J5").Value = 1306
$ws.Range("L5").Value = 2.025592565536499
$ws.Range("M5").Value = 12.80244731903076
$ws.Range("N5").Value = 0.361945199579535
$ws.Range("O5").Value = 0.7331534948783407
$ws.Range("P5").Value = 748
$ws.Range("Q5").Value = 748
$ws.Range("J6").Value = 2443
$ws.Range("L6").Value = 1.560533046722412
$ws.Range("M6").Value = 14.58587741851807
$ws.Range("N6").Value = 0.9983744478347512
$ws.Range("O6").Value = 1.557996304478483
$ws.Range("B7").Value = 5552
$ws.Range("D7").Value = 12.76515109110311
$ws.Range("E7").Value = 12.51972977569507
$ws.Range("F7").Value = 0.9807741158990991
$ws.Range("G7").Value = 2.547037365561873
$ws.Range("H7").Value = -12.50846445560455
$ws.Range("J7").Value = 2444
$ws.Range("L7").Value = 3.443131446838379
$ws.Range("M7").Value = 9.819564819335938
$ws.Range("N7").Value = 0.9842810911257656
$ws.Range("O7").Value = 3.389009187001514
$ws.Range("B8").Value = 5552
$ws.Range("D8").Value = 12.76515109110311
$ws.Range("E8").Value = 12.51972977569507
$ws.Range("F8").Value = 0.9807741158990991
$ws.Range("G8").Value = 2.547037365561873
$ws.Range("H8").Value = -12.50846445560455
$ws.Range("J8").Value = 2434
$ws.Range("L8").Value = 4.242428779602051
$ws.Range("M8").Value = 8.846302032470703
$ws.Range("N8").Value = 0.4402401306434588
$ws.Range("O8").Value = 1.867687381807795
$ws.Range("P8").Value = 423
$ws.Range("Q8").Value = 423
$ws.Range("R8").Value = 423
$ws.Range("B9").Value = 5552
$ws.Range("D9").Value = 12.76515109110311
$ws.Range("E9").Value = 12.51972977569507
$ws.Range("F9").Value = 0.9807741158990991
$ws.Range("G9").Value = 2.547037365561873
$ws.Range("H9").Value = -12.50846445560455
$ws.Range("J9").Value = 1469
$ws.Range("L9").Value = 9.53480339050293
$ws.Range("M9").Value = 5.90083122253418
$ws.Range("N9").Value = 0.4134517271765727
$ws.Range("O9").Value = 3.942180704804922
$ws.Range("P9").Value = 780
$ws.Range("Q9").Value = 1105
$ws.Range("R9").Value = 780
$ws.Range("B10").Value = 5552
$ws.Range("D10").Value = 12.76515109110311
$ws.Range("E10").Value = 12.51972977569507
$ws.Range("F10").Value = 0.9807741158990991
$ws.Range("G10").Value = 2.547037365561873
$ws.Range("H10").Value = -12.50846445560455
$ws.Range("J10").Value = 1309
$ws.Range("L10").Value = 4.989809513092041
$ws.Range("M10").Value = 8.156936645507812
$ws.Range("N10").Value = 0.3532913110586681
$ws.Range("O10").Value = 1.762856248963099
$ws.Range("P10").Value = 130
$ws.Range("Q10").Value = 98
$ws.Range("R10").Value = 98

# --- Cells that must become blank (matching the source data's missing values) ---
$ws.Range("P6").ClearContents()
$ws.Range("Q6").ClearContents()
$ws.Range("R6").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
